$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "test18" environment references to "test21"
$ws.Range("A2").Value = "https://test21.cliotest.com/backoffice/control/main"
$ws.Range("C2").Value = "https://test21.cliotest.com/cabicentral/control/main"
$ws.Range("D2").Value = "https://test21.cliotest.com/warehouse/control/main"
$ws.Range("F2").Value = "virtual_cabitest21"
$ws.Range("G2").Value = "test21"
$ws.Range("K2").Value = "test21"

# Move the active selection from C12 to C13
$ws.Range("C13").Select()
